# Applies the "Updated symbol list" data refresh to Sheet1 of cryptos.xlsx.
# For every data row (2-51) the "Hora" (hour) column G goes from 13 -> 14,
# most "Price" values in column D are refreshed, a couple of "Volume(1h)"
# labels in column E change, and rows 42/43 (BKEXToken <-> CEJI) swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 248.79
$ws.Range("G2").Value = 14
# Row 3
$ws.Range("G3").Value = 14
# Row 4
$ws.Range("D4").Value = 5.449
$ws.Range("G4").Value = 14
# Row 5
$ws.Range("D5").Value = 0.05676
$ws.Range("G5").Value = 14
# Row 6
$ws.Range("D6").Value = 3.386
$ws.Range("G6").Value = 14
# Row 7
$ws.Range("D7").Value = 0.8051
$ws.Range("G7").Value = 14
# Row 8
$ws.Range("G8").Value = 14
# Row 9
$ws.Range("D9").Value = 0.1459
$ws.Range("G9").Value = 14
# Row 10
$ws.Range("D10").Value = 0.07713
$ws.Range("G10").Value = 14
# Row 11
$ws.Range("D11").Value = 0.03177
$ws.Range("G11").Value = 14
# Row 12
$ws.Range("D12").Value = 0.03048
$ws.Range("G12").Value = 14
# Row 13
$ws.Range("D13").Value = 0.09261
$ws.Range("G13").Value = 14
# Row 14
$ws.Range("D14").Value = 3.531
$ws.Range("G14").Value = 14
# Row 15
$ws.Range("D15").Value = 0.001644
$ws.Range("G15").Value = 14
# Row 16
$ws.Range("D16").Value = 0.04715
$ws.Range("G16").Value = 14
# Row 17
$ws.Range("D17").Value = 0.01155
$ws.Range("E17").Value = "16OneONEBestin24h"
$ws.Range("G17").Value = 14
# Row 18
$ws.Range("D18").Value = 0.006362
$ws.Range("G18").Value = 14
# Row 19
$ws.Range("D19").Value = 0.005027
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("G19").Value = 14
# Row 20
$ws.Range("G20").Value = 14
# Row 21
$ws.Range("D21").Value = 0.0001501
$ws.Range("G21").Value = 14
# Row 22
$ws.Range("D22").Value = 0.0003204
$ws.Range("G22").Value = 14
# Row 23
$ws.Range("G23").Value = 14
# Row 24
$ws.Range("D24").Value = 6.425
$ws.Range("G24").Value = 14
# Row 25
$ws.Range("D25").Value = 2.170
$ws.Range("G25").Value = 14
# Row 26
$ws.Range("D26").Value = 0.3304
$ws.Range("G26").Value = 14
# Row 27
$ws.Range("G27").Value = 14
# Row 28
$ws.Range("G28").Value = 14
# Row 29
$ws.Range("G29").Value = 14
# Row 30
$ws.Range("G30").Value = 14
# Row 31
$ws.Range("G31").Value = 14
# Row 32
$ws.Range("G32").Value = 14
# Row 33
$ws.Range("G33").Value = 14
# Row 34
$ws.Range("G34").Value = 14
# Row 35
$ws.Range("G35").Value = 14
# Row 36
$ws.Range("G36").Value = 14
# Row 37
$ws.Range("G37").Value = 14
# Row 38
$ws.Range("G38").Value = 14
# Row 39
$ws.Range("G39").Value = 14
# Row 40
$ws.Range("D40").Value = 0.04074
$ws.Range("G40").Value = 14
# Row 41
$ws.Range("D41").Value = 0.006960
$ws.Range("G41").Value = 14
# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = 0.003504
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = 14
# Row 43
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = 0.1040
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = 14
# Row 44
$ws.Range("D44").Value = 0.007876
$ws.Range("G44").Value = 14
# Row 45
$ws.Range("D45").Value = 0.00005915
$ws.Range("G45").Value = 14
# Row 46
$ws.Range("G46").Value = 14
# Row 47
$ws.Range("D47").Value = 0.0005507
$ws.Range("G47").Value = 14
# Row 48
$ws.Range("G48").Value = 14
# Row 49
$ws.Range("D49").Value = 0.008864
$ws.Range("G49").Value = 14
# Row 50
$ws.Range("G50").Value = 14
# Row 51
$ws.Range("G51").Value = 14
